$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 534.10345
$ws.Range("I28").Value = 155.35
$ws.Range("J28").Value = 1375.7778
$ws.Range("K28").Value = 155.35
$ws.Range("L28").Value = 1375.7778
$ws.Range("M28").Value = 329.65
$ws.Range("N28").Value = -2345.7778
$ws.Range("H69").Value = 23193880
$ws.Range("I69").Value = 3300
$ws.Range("J69").Value = 24224572
$ws.Range("K69").Value = 9900
$ws.Range("L69").Value = 72673716
$ws.Range("M69").Value = -9026
$ws.Range("N69").Value = -72675464
$ws.Range("H70").Value = 4995.727
$ws.Range("I70").Value = 6093.75
$ws.Range("J70").Value = 2067.6667
$ws.Range("K70").Value = 18281.25
$ws.Range("L70").Value = 6203.000100000001
$ws.Range("M70").Value = -18011.25
$ws.Range("N70").Value = -6743.000100000001
$ws.Range("H72").Value = 23193880
$ws.Range("I72").Value = 3300
$ws.Range("J72").Value = 24224572
$ws.Range("K72").Value = 29700
$ws.Range("L72").Value = 218021148
$ws.Range("M72").Value = -25332
$ws.Range("N72").Value = -218029884
$ws.Range("H73").Value = 4995.727
$ws.Range("I73").Value = 6093.75
$ws.Range("J73").Value = 2067.6667
$ws.Range("K73").Value = 18281.25
$ws.Range("L73").Value = 6203.000100000001
$ws.Range("M73").Value = -17345.25
$ws.Range("N73").Value = -8075.000100000001
$ws.Range("H76").Value = 3239.451
$ws.Range("I76").Value = 3028.5862
$ws.Range("J76").Value = 3517.4092
$ws.Range("K76").Value = 3028.5862
$ws.Range("L76").Value = 3517.4092
$ws.Range("M76").Value = -2713.5862
$ws.Range("N76").Value = -4147.4092
$ws.Range("H79").Value = 3239.451
$ws.Range("I79").Value = 3028.5862
$ws.Range("J79").Value = 3517.4092
$ws.Range("K79").Value = 3028.5862
$ws.Range("L79").Value = 3517.4092
$ws.Range("M79").Value = -1936.5862
$ws.Range("N79").Value = -5701.4092
$ws.Range("H98").Value = 1152.2307
$ws.Range("I98").Value = 1043.5454
$ws.Range("J98").Value = 1750
$ws.Range("K98").Value = 1043.5454
$ws.Range("L98").Value = 1750
$ws.Range("M98").Value = 454.4546
$ws.Range("N98").Value = -4746
$ws.Range("H111").Value = 142857660
$ws.Range("I111").Value = 166667230
$ws.Range("J111").Value = 300
$ws.Range("K111").Value = 500001690
$ws.Range("L111").Value = 900
$ws.Range("M111").Value = -499998623
$ws.Range("N111").Value = -7034
$ws.Range("H112").Value = 1447.579
$ws.Range("I112").Value = 460
$ws.Range("J112").Value = 1532.2285
$ws.Range("K112").Value = 1380
$ws.Range("L112").Value = 4596.6855
$ws.Range("M112").Value = -272
$ws.Range("N112").Value = -6812.6855
$ws.Range("H122").Value = 1152.2307
$ws.Range("I122").Value = 1043.5454
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 3130.6362
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -680.6361999999999
$ws.Range("N122").Value = -10150
$ws.Range("H132").Value = 2893.653
$ws.Range("I132").Value = 2088.1162
$ws.Range("J132").Value = 8666.666999999999
$ws.Range("K132").Value = 6264.348599999999
$ws.Range("L132").Value = 26000.001
$ws.Range("M132").Value = -3734.348599999999
$ws.Range("N132").Value = -31060.001
$ws.Range("H135").Value = 997.2222
$ws.Range("I135").Value = 871.875
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 7846.875
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -5311.875
$ws.Range("N135").Value = -23070
$ws.Range("H138").Value = 1912.8
$ws.Range("I138").Value = 798.1786
$ws.Range("J138").Value = 2346.264
$ws.Range("K138").Value = 2394.5358
$ws.Range("L138").Value = 7038.792
$ws.Range("M138").Value = 2745.4642
$ws.Range("N138").Value = -17318.792

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2659.0625
$ws.Range("I2").Value = 2612.0908
$ws.Range("J2").Value = 2762.4
$ws.Range("K2").Value = 2612.0908
$ws.Range("L2").Value = 2762.4
$ws.Range("M2").Value = -2499.0908
$ws.Range("N2").Value = -2988.4
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = None
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = None
$ws.Range("N85").ClearContents()
$ws.Range("H116").Value = 2659.0625
$ws.Range("I116").Value = 2612.0908
$ws.Range("J116").Value = 2762.4
$ws.Range("K116").Value = 2612.0908
$ws.Range("L116").Value = 2762.4
$ws.Range("M116").Value = -318.0907999999999
$ws.Range("N116").Value = -7350.4
$ws.Range("H122").Value = 1283.5
$ws.Range("I122").Value = 777.1539
$ws.Range("J122").Value = 2600
$ws.Range("K122").Value = 2331.4617
$ws.Range("L122").Value = 7800
$ws.Range("M122").Value = 118.5383000000002
$ws.Range("N122").Value = -12700
$ws.Range("H130").Value = 53983
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 53983
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 53983
$ws.Range("N130").Value = -64023
$ws.Range("H132").Value = 1504.675
$ws.Range("I132").Value = 835.1111
$ws.Range("J132").Value = 2895.3076
$ws.Range("K132").Value = 2505.3333
$ws.Range("L132").Value = 8685.9228
$ws.Range("M132").Value = 24.66670000000022
$ws.Range("N132").Value = -13745.9228

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2659.0625
$ws.Range("I3").Value = 2612.0908
$ws.Range("J3").Value = 2762.4
$ws.Range("K3").Value = 2612.0908
$ws.Range("L3").Value = 2762.4
$ws.Range("M3").Value = -2498.0908
$ws.Range("N3").Value = -2990.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2801.8333
$ws.Range("I16").Value = 1505.5
$ws.Range("J16").Value = 3450
$ws.Range("K16").Value = 1505.5
$ws.Range("L16").Value = 3450
$ws.Range("M16").Value = -1218.5
$ws.Range("N16").Value = -4024
$ws.Range("H59").Value = 9715.462
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 9715.462
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 9715.462
$ws.Range("N59").Value = -12005.462
$ws.Range("H113").Value = 2801.8333
$ws.Range("I113").Value = 1505.5
$ws.Range("J113").Value = 3450
$ws.Range("K113").Value = 1505.5
$ws.Range("L113").Value = 3450
$ws.Range("M113").Value = 664.5
$ws.Range("N113").Value = -7790
$ws.Range("H122").Value = 1635.9286
$ws.Range("I122").Value = 1148.625
$ws.Range("J122").Value = 2285.6667
$ws.Range("K122").Value = 3445.875
$ws.Range("L122").Value = 6857.000100000001
$ws.Range("M122").Value = -995.875
$ws.Range("N122").Value = -11757.0001
$ws.Range("H134").Value = 2467.366
$ws.Range("I134").Value = 1473.4445
$ws.Range("J134").Value = 4384.2144
$ws.Range("K134").Value = 4420.333500000001
$ws.Range("L134").Value = 13152.6432
$ws.Range("M134").Value = -1885.333500000001
$ws.Range("N134").Value = -18222.6432

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 4399.385
$ws.Range("I129").Value = 1795.1538
$ws.Range("J129").Value = 5701.5
$ws.Range("K129").Value = 5385.4614
$ws.Range("L129").Value = 17104.5
$ws.Range("M129").Value = -385.4614000000001
$ws.Range("N129").Value = -27104.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1526.7778
$ws.Range("I102").Value = 1526.7778
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1526.7778
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = None
$ws.Range("N102").ClearContents()
$ws.Range("H126").Value = 11112777
$ws.Range("I126").Value = 8335108.5
$ws.Range("J126").Value = 27778784
$ws.Range("K126").Value = 25005325.5
$ws.Range("L126").Value = 83336352
$ws.Range("M126").Value = -25002855.5
$ws.Range("N126").Value = -83341292

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6994783
$ws.Range("I7").Value = 1735.7273
$ws.Range("J7").Value = 45456544
$ws.Range("K7").Value = 1735.7273
$ws.Range("L7").Value = 45456544
$ws.Range("M7").Value = -1623.7273
$ws.Range("N7").Value = -45456768
$ws.Range("H40").Value = 3664.2222
$ws.Range("I40").Value = 3255.6667
$ws.Range("J40").Value = 4174.9165
$ws.Range("K40").Value = 3255.6667
$ws.Range("L40").Value = 4174.9165
$ws.Range("M40").Value = -3119.6667
$ws.Range("N40").Value = -4446.9165
$ws.Range("H122").Value = 1953.3636
$ws.Range("I122").Value = 1973.375
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 5920.125
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -3470.125
$ws.Range("N122").Value = -10600
$ws.Range("H126").Value = 6994783
$ws.Range("I126").Value = 1735.7273
$ws.Range("J126").Value = 45456544
$ws.Range("K126").Value = 5207.1819
$ws.Range("L126").Value = 136369632
$ws.Range("M126").Value = -2737.1819
$ws.Range("N126").Value = -136374572

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2767.1538
$ws.Range("I122").Value = 2533.8333
$ws.Range("J122").Value = 2967.1428
$ws.Range("K122").Value = 7601.499899999999
$ws.Range("L122").Value = 8901.428400000001
$ws.Range("M122").Value = -5151.499899999999
$ws.Range("N122").Value = -13801.4284
$ws.Range("H126").Value = 2467.1936
$ws.Range("I126").Value = 2308.875
$ws.Range("J126").Value = 3010
$ws.Range("K126").Value = 6926.625
$ws.Range("L126").Value = 9030
$ws.Range("M126").Value = -4456.625
$ws.Range("N126").Value = -13970
$ws.Range("H132").Value = 2277.439
$ws.Range("I132").Value = 1366.0526
$ws.Range("J132").Value = 3064.5454
$ws.Range("K132").Value = 4098.1578
$ws.Range("L132").Value = 9193.636200000001
$ws.Range("M132").Value = -1568.1578
$ws.Range("N132").Value = -14253.6362

Write-Output "Applied all profit sheet updates."